$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for updated crypto symbol data.
# Values are written with a leading apostrophe so Excel stores them as literal
# text (matching the source data's text-based price/percentage strings)
# instead of re-interpreting them as numbers or percentages.

$ws.Range("D2").Value = "'308.70"
$ws.Range("E2").Value = "'-0.13%"
$ws.Range("D3").Value = "'39.86"
$ws.Range("E3").Value = "'2.98%"
$ws.Range("D4").Value = "'5.143"
$ws.Range("E4").Value = "'1.08%"
$ws.Range("D5").Value = "'0.08147"
$ws.Range("E5").Value = "'-0.43%"
$ws.Range("D6").Value = "'1.945"
$ws.Range("E6").Value = "'-1.66%"
$ws.Range("E7").Value = "'3.52%"
$ws.Range("D8").Value = "'4.225"
$ws.Range("E8").Value = "'0.91%"
$ws.Range("D9").Value = "'0.9299"
$ws.Range("E9").Value = "'-0.22%"
$ws.Range("D10").Value = "'0.1442"
$ws.Range("E10").Value = "'2.89%"
$ws.Range("D11").Value = "'0.1922"
$ws.Range("E11").Value = "'-1.45%"
$ws.Range("D12").Value = "'0.09095"
$ws.Range("E12").Value = "'-2.14%"
$ws.Range("D13").Value = "'0.03525"
$ws.Range("E13").Value = "'2.48%"
$ws.Range("D14").Value = "'0.09788"
$ws.Range("E14").Value = "'-0.67%"
$ws.Range("D15").Value = "'0.001397"
$ws.Range("E15").Value = "'-0.92%"
$ws.Range("D16").Value = "'0.005908"
$ws.Range("E16").Value = "'-3.10%"
$ws.Range("D17").Value = "'3.922"
$ws.Range("E17").Value = "'6.98%"
$ws.Range("D18").Value = "'3.346"
$ws.Range("E18").Value = "'-4.09%"
$ws.Range("D19").Value = "'0.3430"
$ws.Range("E19").Value = "'-0.57%"
$ws.Range("D20").Value = "'0.1311"
$ws.Range("E20").Value = "'-1.43%"
$ws.Range("D21").Value = "'4.694"
$ws.Range("E21").Value = "'-2.45%"
$ws.Range("D22").Value = "'0.2425"
$ws.Range("E22").Value = "'-1.18%"
$ws.Range("D23").Value = "'0.04382"
$ws.Range("E23").Value = "'-1.77%"
$ws.Range("D24").Value = "'0.001229"
$ws.Range("E24").Value = "'-0.75%"
$ws.Range("D25").Value = "'0.004376"
$ws.Range("E25").Value = "'4.84%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'-0.08%"
$ws.Range("D27").Value = "'0.0004002"
$ws.Range("E27").Value = "'-10.01%"
$ws.Range("D39").Value = "'0.02057"
$ws.Range("E39").Value = "'-2.87%"
$ws.Range("D40").Value = "'0.05070"
$ws.Range("E40").Value = "'-1.94%"
$ws.Range("D41").Value = "'0.007419"
$ws.Range("E41").Value = "'-0.74%"
$ws.Range("D42").Value = "'0.009884"
$ws.Range("E42").Value = "'-1.23%"
$ws.Range("E43").Value = "'-0.28%"
$ws.Range("D44").Value = "'0.002131"
$ws.Range("E44").Value = "'-0.08%"
$ws.Range("D45").Value = "'0.009376"
$ws.Range("E45").Value = "'-3.07%"
$ws.Range("D46").Value = "'0.00006395"
$ws.Range("E46").Value = "'1.33%"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("D48").Value = "'0.002715"
$ws.Range("E49").Value = "'-18.83%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.08%"
